{"js": "const body = context.document.body;\nconst replacements = [\n  [\"2026-01-19 Monday\", \"2026-01-20 Tuesday\"],\n  [\"73\u00d732=2336\", \"34\u00d765=2210\"],\n  [\"66\u00d744=2904\", \"26\u00d746=1196\"],\n  [\"84\u00d757=4788\", \"26\u00d754=1404\"],\n  [\"80\u00d761=4880\", \"74\u00d794=6956\"],\n  [\"52\u00d771=3692\", \"79\u00d773=5767\"],\n  [\"35\u00d717=595\", \"60\u00d726=1560\"],\n  [\"91\u00d716=1456\", \"35\u00d796=3360\"],\n  [\"42\u00d711=462\", \"11\u00d790=990\"],\n  [\"14\u00d799=1386\", \"44\u00d777=3388\"],\n  [\"55\u00d712=660\", \"47\u00d774=3478\"],\n  [\"77\u00d711=847\", \"65\u00d752=3380\"],\n  [\"63\u00d762=3906\", \"46\u00d729=1334\"],\n  [\"64\u00d761=3904\", \"90\u00d733=2970\"],\n  [\"65\u00d768=4420\", \"46\u00d721=966\"],\n  [\"66\u00d757=3762\", \"29\u00d720=580\"],\n  [\"94\u00d726=2444\", \"56\u00d756=3136\"],\n  [\"55\u00d723=1265\", \"77\u00d761=4697\"],\n  [\"31\u00d737=1147\", \"63\u00d737=2331\"],\n  [\"78\u00d711=858\", \"46\u00d727=1242\"],\n  [\"88\u00d788=7744\", \"89\u00d714=1246\"],\n  [\"19\u00d737=703\", \"71\u00d787=6177\"],\n  [\"91\u00d726=2366\", \"74\u00d740=2960\"],\n  [\"85\u00d751=4335\", \"49\u00d751=2499\"],\n  [\"89\u00d770=6230\", \"45\u00d745=2025\"],\n  [\"40\u00d761=2440\", \"62\u00d750=3100\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Not found: ' + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nreturn 'done: ' + replacements.length + ' replacements';", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-19 Monday\", \"2026-01-20 Tuesday\"),\n    @(\"73\u00d732=2336\", \"34\u00d765=2210\"),\n    @(\"66\u00d744=2904\", \"26\u00d746=1196\"),\n    @(\"84\u00d757=4788\", \"26\u00d754=1404\"),\n    @(\"80\u00d761=4880\", \"74\u00d794=6956\"),\n    @(\"52\u00d771=3692\", \"79\u00d773=5767\"),\n    @(\"35\u00d717=595\", \"60\u00d726=1560\"),\n    @(\"91\u00d716=1456\", \"35\u00d796=3360\"),\n    @(\"42\u00d711=462\", \"11\u00d790=990\"),\n    @(\"14\u00d799=1386\", \"44\u00d777=3388\"),\n    @(\"55\u00d712=660\", \"47\u00d774=3478\"),\n    @(\"77\u00d711=847\", \"65\u00d752=3380\"),\n    @(\"63\u00d762=3906\", \"46\u00d729=1334\"),\n    @(\"64\u00d761=3904\", \"90\u00d733=2970\"),\n    @(\"65\u00d768=4420\", \"46\u00d721=966\"),\n    @(\"66\u00d757=3762\", \"29\u00d720=580\"),\n    @(\"94\u00d726=2444\", \"56\u00d756=3136\"),\n    @(\"55\u00d723=1265\", \"77\u00d761=4697\"),\n    @(\"31\u00d737=1147\", \"63\u00d737=2331\"),\n    @(\"78\u00d711=858\", \"46\u00d727=1242\"),\n    @(\"88\u00d788=7744\", \"89\u00d714=1246\"),\n    @(\"19\u00d737=703\", \"71\u00d787=6177\"),\n    @(\"91\u00d726=2366\", \"74\u00d740=2960\"),\n    @(\"85\u00d751=4335\", \"49\u00d751=2499\"),\n    @(\"89\u00d770=6230\", \"45\u00d745=2025\"),\n    @(\"40\u00d761=2440\", \"62\u00d750=3100\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}"}
